# locator_master.xlsx update: add "Checkout" title locators for the
# checkout page and the payment page on both the android and ios sheets.
#
# android (sheet1) uses ID locators pointing at view ids; ios (sheet2) uses
# ACCESSIBILITY_ID locators pointing at the "Checkout" label text.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # android
$ws2 = $wb.Worksheets.Item(2)   # ios

# --- Make room for the two new locator rows on each sheet --------------
# android: new row before "checkout_page_name" (old row 14) and before
# "payment_page_card_name" (old row 20, becomes row 22 after first insert).
$ws1.Rows.Item(14).Insert()
$ws1.Rows.Item(21).Insert()

# ios: new row before "my_cart_page_checkout_button" (old row 10) and
# before "payment_page_card_name" (old row 20, becomes row 22 after the
# first insert).
$ws2.Rows.Item(10).Insert()
$ws2.Rows.Item(21).Insert()

# --- Fill in the new cells ----------------------------------------------
# Order chosen so new shared strings are created in this sequence:
# "Checkout", "checkout_page_checkout_title",
# "com.saucelabs.mydemoapp.android:id/checkoutTitleTV",
# "payment_page_checkout_title",
# "com.saucelabs.mydemoapp.android:id/enterPaymentTitleTV"
$ws2.Range("C10").Value = "Checkout"
$ws2.Range("A10").Value = "checkout_page_checkout_title"
$ws1.Range("C14").Value = "com.saucelabs.mydemoapp.android:id/checkoutTitleTV"
$ws2.Range("A21").Value = "payment_page_checkout_title"
$ws1.Range("C21").Value = "com.saucelabs.mydemoapp.android:id/enterPaymentTitleTV"

$ws1.Range("A14").Value = "checkout_page_checkout_title"
$ws1.Range("B14").Value = "ID"

$ws1.Range("A21").Value = "payment_page_checkout_title"
$ws1.Range("B21").Value = "ID"

$ws2.Range("B10").Value = "ACCESSIBILITY_ID"

$ws2.Range("B21").Value = "ACCESSIBILITY_ID"
$ws2.Range("C21").Value = "Checkout"

# --- View state: user ended up on the ios sheet, with android's selection
#     left at F24 and ios's selection at E18 -----------------------------
$ws1.Range("F24").Select()
$ws2.Activate()
$ws2.Range("E18").Select()
